# Knights of the Force Academy - language.xlsx localization sheet update
# "more art assets for act 3, added simpler first level, added tutorial"
#
# - newton_first_law_dlg_6 text is reworded.
# - four new tutorial-help rows (dragForcePushHelp, dragForceAdjustHelp,
#   playButtonHelp, trashButtonHelp) are appended after the existing data.
# - the sheet view is scrolled / re-selected to show the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Act 3 tutorial help strings (rows 95-98).
$ws.Range("A95").Value = "dragForcePushHelp"
$ws.Range("A96").Value = "dragForceAdjustHelp"
$ws.Range("B96").Value = "You can reposition a Force Field by dragging it."
$ws.Range("B95").Value = "Drag the Force Field to this surface."
$ws.Range("A97").Value = "playButtonHelp"
$ws.Range("B97").Value = "Press this button to start the simulation."
$ws.Range("A98").Value = "trashButtonHelp"

# Reword the existing "Gather your forces..." dialogue line (key unchanged).
$ws.Range("B47").Value = "Let us summon our mighty knights for a daring rescue!"

$ws.Range("B98").Value = "Press and hold the trash button to clear out the Force Fields."

$ws.Range("C96").Value = 3
$ws.Range("C97").Value = 3

# Scroll the sheet view down to the newly added rows and select B98.
$win = $excel.ActiveWindow
$win.ScrollRow = 79
$win.ScrollColumn = 1
$ws.Range("B98").Select()
